$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.138.66"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.912.08"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.45"
$ws.Range("E5").Value = "  -3.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.07"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.910.22"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.97"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.60"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.396.51"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.095.86"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.915.44"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.14"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.653"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.94"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.00"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.27"
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("E29").Value = "  +7.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -4.95%  "
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.01"
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.69"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.955"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.40"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.87"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.91"
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("E40").Value = "  -4.81%  "
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.16"
$ws.Range("E42").Value = "  +5.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.17"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.268"
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.731.25"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.57"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0337"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "348.97"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000212"
$ws.Range("E51").Value = "  +9.74%  "

Write-Output "Applied cryptos update"
